# Homework 2 update: add a second worksheet ("Sheet2") with flow/levee data,
# plus tweak the selection state that was left on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: just the view/selection changed (tabSelected moves to new sheet) ---
$ws1.Range("C1:Q18").Select() | Out-Null

# --- Add Sheet2 right after Sheet1 (becomes the active tab) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Row 2 headers
$ws2.Range("B2").Value = "Stage 200"
$ws2.Range("F2").Value = "Choose to build this many"
$ws2.Range("Q2").Value = "Choice"

# Row 3 headers
$ws2.Range("B3").Value = "Have this many"
$ws2.Range("C3").Value = "Height"
$ws2.Range("D3").Value = "Mean Flow"
$ws2.Range("E3").Value = "Flow SD"

$row3vals = @(0,1,2,3,4,5,6,7,8,9,10)
$col = 6
foreach ($v in $row3vals) {
    $ws2.Cells.Item(3, $col).Value = $v
    $col++
}
$ws2.Range("Q3").Value = "deltaH*"
$ws2.Range("R3").Value = "f*"

# Row 4
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 10
$ws2.Range("E4").Value = -1
$row4vals = @(999,900,800,650,665,700,750,800,900,1000,1010)
$col = 6
foreach ($v in $row4vals) {
    $ws2.Cells.Item(4, $col).Value = $v
    $col++
}
$ws2.Range("Q4").Value = 3
$ws2.Range("R4").Value = 650

# Row 5
$ws2.Range("C5").Value = 0
$ws2.Range("D5").Value = 20
$ws2.Range("E5").Value = -1
$ws2.Range("F5").Formula = "=F4+100"
$ws2.Range("G5:P5").Formula = "=G4+100"
$ws2.Range("Q5").Value = 3
$ws2.Range("R5").Formula = "=MIN(F5:P5)"

# Row 6-8 (no height/flow data, just id columns + Choice)
$ws2.Range("C6").Value = 0
$ws2.Range("D6").Value = 30
$ws2.Range("E6").Value = -1
$ws2.Range("Q6").Value = 2

$ws2.Range("C7").Value = 0
$ws2.Range("D7").Value = 40
$ws2.Range("E7").Value = -1
$ws2.Range("Q7").Value = 2

$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 50
$ws2.Range("E8").Value = -1
$ws2.Range("Q8").Value = 1

# Row 9
$ws2.Range("C9").Value = 0
$ws2.Range("D9").Value = 10
$ws2.Range("E9").Value = 0
$ws2.Range("Q9").Value = 1

# Row 10
$ws2.Range("C10").Value = 0
$ws2.Range("D10").Value = 20
$ws2.Range("E10").Value = 0
$ws2.Range("P10").Style = "Bad"
$ws2.Range("Q10").Value = 1

# Row 11
$ws2.Range("C11").Value = 0
$ws2.Range("D11").Value = 30
$ws2.Range("E11").Value = 0
$ws2.Range("O11:P11").Style = "Bad"
$ws2.Range("Q11").Value = 1

# Row 12
$ws2.Range("C12").Value = 0
$ws2.Range("D12").Value = 40
$ws2.Range("E12").Value = 0
$ws2.Range("N12:P12").Style = "Bad"
$ws2.Range("Q12").Value = 1

# Row 13
$ws2.Range("C13").Value = 0
$ws2.Range("D13").Value = 50
$ws2.Range("E13").Value = 0
$ws2.Range("M13:P13").Style = "Bad"
$ws2.Range("Q13").Value = 0

# Row 14
$ws2.Range("C14").Value = 0
$ws2.Range("D14").Value = 10
$ws2.Range("E14").Value = 1
$ws2.Range("L14:P14").Style = "Bad"
$ws2.Range("Q14").Value = 0

# Row 15
$ws2.Range("C15").Value = 0
$ws2.Range("D15").Value = 20
$ws2.Range("E15").Value = 1
$ws2.Range("K15:P15").Style = "Bad"
$ws2.Range("Q15").Value = 0

# Row 16
$ws2.Range("C16").Value = 0
$ws2.Range("D16").Value = 30
$ws2.Range("E16").Value = 1
$ws2.Range("J16:P16").Style = "Bad"
$ws2.Range("Q16").Value = 0

# Row 17
$ws2.Range("C17").Value = 0
$ws2.Range("D17").Value = 40
$ws2.Range("E17").Value = 1
$ws2.Range("I17:P17").Style = "Bad"
$ws2.Range("Q17").Value = 0

# Row 18
$ws2.Range("C18").Value = 0
$ws2.Range("D18").Value = 50
$ws2.Range("E18").Value = 1
$ws2.Range("H18:P18").Style = "Bad"
$ws2.Range("Q18").Value = 0

# Row 19
$ws2.Range("C19").Value = 1
$ws2.Range("G19:P19").Style = "Bad"
$ws2.Range("Q19").Value = 0
$ws2.Range("R19").Value = "?"

# Rows 20-33: just an incrementing id column
$idvals = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15)
$r = 20
foreach ($v in $idvals) {
    $ws2.Cells.Item($r, 3).Value = $v
    $r++
}

# Column widths for D (Mean Flow) and E (Flow SD) so labels fit
# (closest achievable values to Excel's original bestFit widths of
# 10.7109375 / 13.5703125 given this engine's column-width quantization)
$ws2.Columns.Item(4).ColumnWidth = 9.833333333333332
$ws2.Columns.Item(5).ColumnWidth = 12.666666666666666

# Sheet2 view/selection
$ws2.Range("H10").Select() | Out-Null
